{"js": "// Fill in the previously-empty \"value\" cells for the RETENTION ratio and\n// the Answer Recall (ARL/ARS/ARA) rows with their computed results.\n// The target document has a single table; each row we care about has a\n// label in its first cell and an (until now empty) value in its last cell.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.rows.load(\"items\");\nawait context.sync();\n\nconst rows = table.rows.items;\nfor (const row of rows) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\nfor (const row of rows) {\n  for (const cell of row.cells.items) {\n    cell.body.load(\"text\");\n  }\n}\nawait context.sync();\n\n// label (as it appears in the first cell of the row) -> computed value to\n// write into the row's last (still empty) cell.\nconst updates = [\n  [\"Ratio\", \"0.8\"],\n  [\"Answer Recall Lenient (ARL)\", \"0.4285\"],\n  [\"Answer Recall Strict (ARS)\", \"0.2857\"],\n  [\"Answer Recall Average (ARA)\", \"0.3571\"],\n];\n\nfor (const [label, value] of updates) {\n  const row = rows.find(\n    (r) => r.cells.items.length > 1 && r.cells.items[0].body.text.trim() === label\n  );\n  if (!row) {\n    throw new Error(`Could not find row labelled \"${label}\"`);\n  }\n  const cells = row.cells.items;\n  const valueCell = cells[cells.length - 1];\n  const para = valueCell.body.paragraphs.getFirst();\n  const range = para.insertText(value, \"Replace\");\n  range.font.bold = true;\n  range.font.size = 12;\n  range.font.sizeBidirectional = 12;\n}\n\nawait context.sync();\n", "ps1": "# Fill in the previously-empty \"value\" cells for the RETENTION ratio and\n# the Answer Recall (ARL/ARS/ARA) rows with their computed results.\n# The document has a single table; the rows we care about carry a label\n# in their first cell and an (until now empty) value in their last cell.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$rowCount = $t.Rows.Count\n\n# label (as it appears in the first cell of the row) -> computed value to\n# write into the row's last (still empty) cell.\n$updates = [ordered]@{\n    \"Ratio\"                        = \"0.8\"\n    \"Answer Recall Lenient (ARL)\"  = \"0.4285\"\n    \"Answer Recall Strict (ARS)\"   = \"0.2857\"\n    \"Answer Recall Average (ARA)\"  = \"0.3571\"\n}\n\nforeach ($label in $updates.Keys) {\n    $value = $updates[$label]\n    $matched = $false\n    for ($i = 1; $i -le $rowCount; $i++) {\n        $row = $t.Rows.Item($i)\n        $cellCount = $row.Cells.Count\n        if ($cellCount -lt 2) { continue }\n        # Cell text ends with a cell-mark (CR + BEL); strip those before comparing.\n        $firstText = ($row.Cells.Item(1).Range.Text -replace \"[\\r\\a]\", \"\").Trim()\n        if ($firstText -eq $label) {\n            $lastCell = $row.Cells.Item($cellCount)\n            $lastCell.Range.Text = $value\n            $lastCell.Range.Font.Bold = 1\n            $lastCell.Range.Font.Size = 12\n            $lastCell.Range.Font.SizeBi = 12\n            $matched = $true\n            break\n        }\n    }\n    if (-not $matched) {\n        Write-Output \"NOT FOUND: $label\"\n    }\n}\n"}
